$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45973
$ws.Range("B2").Value = 36.9
$ws.Range("C2").Value = 18.71
$ws.Range("D2").Value = 15.77
$ws.Range("E2").Value = 12.26
$ws.Range("F2").Value = 9.210000000000001
$ws.Range("G2").Value = 14.48
$ws.Range("H2").Value = 37.57
$ws.Range("I2").Value = 52.36
$ws.Range("J2").Value = 45.46
$ws.Range("K2").Value = 11.59
$ws.Range("L2").Value = 2.09
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 1.34
$ws.Range("O2").Value = 0.79
$ws.Range("P2").Value = 1.15
$ws.Range("Q2").Value = 3.24
$ws.Range("R2").Value = 16.15
$ws.Range("S2").Value = 46.42
$ws.Range("T2").Value = 55.23
$ws.Range("U2").Value = 63.5
$ws.Range("V2").Value = 71.56999999999999
$ws.Range("W2").Value = 65.66
$ws.Range("X2").Value = 51.95
$ws.Range("Y2").Value = 33.36
$ws.Range("Z2").Value = 27.78
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 55.64
$ws.Range("AD2").Value = 68.61
$ws.Range("AF2").Value = 59.36
$ws.Range("AG2").Value = "1h-16h"
